$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New columns J (unique non-blank part numbers) and K (count of rows with
# that part number that still need ordering, i.e. no free-stock note in H)
# ---------------------------------------------------------------------------

# Header cells
$ws.Range("J1").Value = "uniq"
$ws.Range("K1").Value = "cnt non free"
$ws.Range("J1:K1").Font.Bold = $true

# J2:J49 - array formula that walks through G2:G52 picking out each new
# unique value in turn (classic "extract unique values" INDEX/MATCH/COUNTIF
# array formula idiom).
for ($r = 2; $r -le 49; $r++) {
    $prevRow = $r - 1
    $formula = '=INDEX($G$2:$G$52, MATCH(0, COUNTIF($J$1:J' + $prevRow + ', $G$2:$G$52), 0))'
    $ws.Range("J$r").FormulaArray = $formula
}

# K2:K31 - count, for each unique part in J, how many rows reference it and
# have an empty "PCB Train Free" column H (i.e. still need to be ordered).
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("K$r").Formula = '=COUNTIFS($G$2:$G$75,J' + $r + ',$H$2:$H$75,"")'
}

# Formatting for the J column values: Consolas 10, dark grey, left/centered.
$jvals = $ws.Range("J2:J49")
$jvals.Font.Name = "Consolas"
$jvals.Font.Size = 10
$jvals.Font.Color = 2250274
$jvals.HorizontalAlignment = -4131
$jvals.VerticalAlignment = -4108

# Conditional formatting on column K: highlight zero counts (nothing left to
# order) by fading the text to near-white.
$cf = $ws.Range("K1:K1048576").FormatConditions.Add(1, 3, "0")
$cf.Font.ColorIndex = 1
$cf.Font.TintAndShade = -0.14996795556505021

# ---------------------------------------------------------------------------
# Misc layout tweaks that came along with this change
# ---------------------------------------------------------------------------

# Column H (free-stock note) needed to be much wider once J/K were added.
$ws.Columns("H").ColumnWidth = 35.7109375

# View: zoomed out a bit further and the last touched cell was K49.
$ws.Application.ActiveWindow.Zoom = 55
$ws.Range("K49").Select()
